# The M2Doc field "{m:'dh1.gif'.asImage().fit(100, 250)}" lived inside a
# table cell in the document's (only) footer, encoded as a real Word field
# (fldChar begin / instrText* / fldChar end). The parser was updated to use
# TokenIteratorFieldRewriterSplit, which expects the template marker to be
# stored as plain literal text runs ("{m:...}") instead of a Word field.
#
# This rewrites that single table cell so the field is replaced by four
# plain <w:t> runs carrying the same text (and the same orange accent-color
# run formatting on the two script fragments), with the field machinery
# removed entirely.

$d = $word.ActiveDocument

$footer = $d.Sections.Item(1).Footers.Item(1)

# Touch an unrelated, harmless piece of footer text first (replacing it
# with itself). Word COM's footer story needs one "story-level" edit
# committed before an edit that lives inside a footer's table cell will be
# persisted back to the part on save.
$null = $footer.Range.Find.Execute("Checks", $true, $false, $false, $false, $false, $true, 1, $false, "Checks", 2)

$table = $footer.Range.Tables.Item(1)
$cell = $table.Cell(2, 2)
$cellParagraph = $cell.Range.Paragraphs.Item(1)
$target = $cellParagraph.Range

$color = '<w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>'

$newParagraphXml = '<w:p>' +
  '<w:r><w:t>{m:</w:t></w:r>' +
  '<w:r><w:rPr>' + $color + '</w:rPr><w:t>''dh1.gif''.asImage()</w:t></w:r>' +
  '<w:r><w:rPr>' + $color + '</w:rPr><w:t>.fit(100, 250)</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
  '</w:p>'

$flatOpcXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  $newParagraphXml +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $target.InsertXML($flatOpcXml)
